$wb = $excel.ActiveWorkbook

# Values (accuracy flags) to write into the newly-inserted column(s),
# one per data row (rows 2-14). Row 1 is the header row and is left blank
# in the new column(s), matching the reference edit.
$vals = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 0
}

# --- Sheet "sheet1": insert TWO columns before D (old D->F, old E->G) ---
$ws1 = $wb.Worksheets.Item("sheet1")
$ws1.Columns("D:E").Insert()
foreach ($r in $vals.Keys) {
    $ws1.Cells.Item($r, 4).Value2 = $vals[$r]
    $ws1.Cells.Item($r, 5).Value2 = $vals[$r]
}

# --- Sheet "Feuille2": insert ONE column before D (old D->E, old E->F) ---
$ws2 = $wb.Worksheets.Item("Feuille2")
$ws2.Columns("D:D").Insert()
foreach ($r in $vals.Keys) {
    $ws2.Cells.Item($r, 4).Value2 = $vals[$r]
}

# --- Sheet "Feuille3": insert ONE column before D (old D->E, old E->F) ---
$ws3 = $wb.Worksheets.Item("Feuille3")
$ws3.Columns("D:D").Insert()
$ws3.Columns("D:D").ClearFormats()
foreach ($r in $vals.Keys) {
    $ws3.Cells.Item($r, 4).Value2 = $vals[$r]
}

# --- Sheet "Feuille4": left untouched ---
